# Atoms.xlsx: rename the "Identifier" attribute-name cell on both the
# Concepts sheet and the Relations sheet to sheet-specific names, and
# update the selected cell / active sheet to reflect where the user
# ended up after making the edits.

$wb = $excel.ActiveWorkbook
$wsConcepts  = $wb.Worksheets.Item("Concepts")
$wsRelations = $wb.Worksheets.Item("Relations")

# Relations sheet: the [Relation] atom's "name" attribute was called
# "Identifier"; rename it to "RelationName".
$wsRelations.Range("B2").Value = "RelationName"

# Concepts sheet: the [Concept] atom's "name" attribute was called
# "Identifier"; rename it to "ConceptName".
$wsConcepts.Range("B2").Value = "ConceptName"

# Leave a selection on the Relations sheet at B5 ...
$wsRelations.Activate()
$wsRelations.Range("B5").Select()

# ... but finish with the Concepts sheet active and selected at D6,
# matching the saved view state of the workbook.
$wsConcepts.Activate()
$wsConcepts.Range("D6").Select()
